$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping-calefaccion-detalle")

$ws.Range("A1").Value = "Sin instalación pero con aparatos que permiten calentar alguna habitación"
$ws.Range("A2").Value = "Con calefacción individual"
$ws.Range("A3").Value = "Sin ningún medio"
$ws.Range("A4").Value = "Con calefacción colectiva"
